# Task Assignment List - 12Dec17 update
# Assign tasks to each team member on the "Rebuild Tasks" sheet, and
# tidy up the row heights / selection left over from that editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rebuild Tasks")

# --- Assign tasks to each team member (column B) ---
$ws.Range("B3").Value = "Jennifer"
$ws.Range("B4").Value = "Jennifer"
$ws.Range("B5").Value = "Andrea"
$ws.Range("B6").Value = "Clark"
$ws.Range("B9").Value = "Alex"
$ws.Range("B10").Value = "Alex"
$ws.Range("B11").Value = "Eliseo"
$ws.Range("B12").Value = "Andrea"
$ws.Range("B13").Value = "Eliseo"

# B13 previously had no horizontal alignment set (unlike the rest of the
# "Assigned to" column); bring it in line with the other rows now that it
# has content.
$ws.Range("B13").HorizontalAlignment = -4108

# --- Row-height tweaks from resizing the header rows while reviewing ---
$ws.Rows("1").RowHeight = 23.25
$ws.Rows("2").RowHeight = 54.75

# --- Leave the selection where the user finished editing ---
$ws.Range("I17").Select() | Out-Null
